$d = $word.ActiveDocument

# Replace all occurrences of "June 26, 2022" with "June 29, 2022"
$d.Content.Find.Execute("June 26, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "June 29, 2022", 2)

# Replace "August 25, 2022" with "August 28, 2022"
$d.Content.Find.Execute("August 25, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "August 28, 2022", 2)
